$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." spread across four runs:
#   "Versi" | "on" | " 2" | "."
# with a (hidden) _GoBack bookmark sitting between the " 2" and "." runs.
# The target reads "Version 1." as two runs ("Version" and " 1.") with the
# _GoBack bookmark moved to sit right after the merged " 1." run.

# 1. Drop the existing _GoBack bookmark first so the text edits below
#    (which span its location) aren't blocked/split by it.
$content = $d.Content
$goBack = $content.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Merge "Versi" + "on" into a single "Version" run.
$word1 = $d.Range(0, 7)
$word1.Find.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# 3. Change "2" to "1" and merge " 2" into a single " 1" run (the trailing
#    "." stays a separate run for now).
$word2 = $d.Content
$word2.Find.Execute(" 2", $true, $false, $false, $false, $false, $true, 1, $false, " 1", 2)

# 4. Re-create the _GoBack bookmark, collapsed right after " 1" (i.e.
#    exactly where it used to sit, between " 1"/" 2" and ".").
$bmSpot = $d.Range(9, 9)
$bmSpot.Bookmarks.Add("_GoBack")

# 5. Remove the now-separate "." run text...
$trailingDot = $d.Range(9, 10)
$trailingDot.Delete()

# 6. ...and re-insert it immediately before the bookmark's (live) range so
#    it merges back into the preceding " 1" run, becoming " 1.", with the
#    bookmark ending up right after it - matching the target structure.
$bmRange = $d.Bookmarks.Item("_GoBack").Range
$bmRange.InsertBefore(".")
